# Atualização automática de DOM_PEDRITO.xlsx
#
# 1. Renomeia a aba "Paineis DARQ" para "PAINEIS DARQ"
# 2. Renomeia a aba "Recolhimento x Eliminacao" para "RECOLHIMENTO X ELIMINAÇÃO"
# 3. Remove a aba "Desarquivamentos Pendentes"

$wb = $excel.ActiveWorkbook

# Avoid the "delete sheet" confirmation prompt Excel normally shows.
$excel.DisplayAlerts = $false

$wb.Worksheets.Item("Paineis DARQ").Name = "PAINEIS DARQ"
$wb.Worksheets.Item("Recolhimento x Eliminacao").Name = "RECOLHIMENTO X ELIMINAÇÃO"
$wb.Worksheets.Item("Desarquivamentos Pendentes").Delete()

$excel.DisplayAlerts = $true
